# Actualización automática 2025-11-25 08:30:08
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M4").Value  = 8453.51
$wsGrupo.Range("L5").Value  = 3184.26
$wsGrupo.Range("M5").Value  = 259.1
$wsGrupo.Range("M12").Value = 3218.23
$wsGrupo.Range("I24").Value = 50.4
$wsGrupo.Range("M24").Value = 122.28
$wsGrupo.Range("H29").Value = 410.4
$wsGrupo.Range("L29").Value = 2721.89
$wsGrupo.Range("M29").Value = 4169.38
$wsGrupo.Range("D36").Value = 2436.48
$wsGrupo.Range("L36").Value = 3158.97
$wsGrupo.Range("H47").Value = 1910.7

$wsGrupo.Range("D56").Value = "8 de 54"
$wsGrupo.Range("H56").Value = "5 de 54"
$wsGrupo.Range("I56").Value = "5 de 54"
$wsGrupo.Range("M56").Value = "17 de 54"

# ---------------------------------------------------------------------
# Sheet: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value  = 10987.55
$wsMensual.Range("F5").Value  = 3443.36
$wsMensual.Range("F12").Value = 3693.43
$wsMensual.Range("F24").Value = 4553.84
$wsMensual.Range("F29").Value = 8902.23
$wsMensual.Range("F36").Value = 5652.31
$wsMensual.Range("F47").Value = 1910.7
$wsMensual.Range("F60").Value = 76700.61

# ---------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column E width grows from 23 to 24 characters
$wsCumplimiento.Columns.Item(5).ColumnWidth = 23.1

$wsCumplimiento.Range("D3").Value  = 6903.65
$wsCumplimiento.Range("E3").Value  = -280.3899999999994
$wsCumplimiento.Range("F3").Value  = 1.042334137569716

$wsCumplimiento.Range("D6").Value  = 4046.4
$wsCumplimiento.Range("E6").Value  = -1138.81631853974
$wsCumplimiento.Range("F6").Value  = 1.39167103798292

$wsCumplimiento.Range("D7").Value  = 518.4
$wsCumplimiento.Range("E7").Value  = 801.6
$wsCumplimiento.Range("F7").Value  = 0.3927272727272727

$wsCumplimiento.Range("D11").Value = 19287.88
$wsCumplimiento.Range("E11").Value = -5051.890000000001
$wsCumplimiento.Range("F11").Value = 1.354867487262916

$wsCumplimiento.Range("D12").Value = 41672.29
$wsCumplimiento.Range("E12").Value = 23271.71
$wsCumplimiento.Range("F12").Value = 0.641664972899729

$wsCumplimiento.Range("D14").Value = 74826.70000000001
$wsCumplimiento.Range("E14").Value = 24129.55685923838
$wsCumplimiento.Range("F14").Value = 0.756159361468555
